$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1208.8
$ws.Range("J17").Value = 1208.8
$ws.Range("L17").Value = 3626.4
$ws.Range("N17").Value = -3962.4
$ws.Range("H53").Value = 248.375
$ws.Range("I53").Value = 147.75
$ws.Range("K53").Value = 147.75
$ws.Range("M53").Value = 489.25
$ws.Range("H62").Value = 4320.2
$ws.Range("I62").Value = 3901
$ws.Range("J62").Value = 4425
$ws.Range("K62").Value = 3901
$ws.Range("L62").Value = 4425
$ws.Range("M62").Value = -3277
$ws.Range("N62").Value = -5673
$ws.Range("H65").Value = 4320.2
$ws.Range("I65").Value = 3901
$ws.Range("J65").Value = 4425
$ws.Range("K65").Value = 19505
$ws.Range("L65").Value = 22125
$ws.Range("M65").Value = -16385
$ws.Range("N65").Value = -28365
$ws.Range("H92").Value = 377.66666
$ws.Range("I92").Value = 296.5
$ws.Range("J92").Value = 540
$ws.Range("K92").Value = 296.5
$ws.Range("L92").Value = 540
$ws.Range("M92").Value = 951.5
$ws.Range("N92").Value = -3036
$ws.Range("H141").Value = 6276.385
$ws.Range("I141").Value = 4836.5625
$ws.Range("J141").Value = 8580.1
$ws.Range("K141").Value = 14509.6875
$ws.Range("L141").Value = 25740.3
$ws.Range("M141").Value = -9329.6875
$ws.Range("N141").Value = -36100.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 3966.6667
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 3966.6667
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 3966.6667
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -4266.6667
$ws.Range("H32").Value = 6141.75
$ws.Range("I32").Value = 4410.4165
$ws.Range("J32").Value = 15231.25
$ws.Range("K32").Value = 4410.4165
$ws.Range("L32").Value = 15231.25
$ws.Range("M32").Value = -4123.4165
$ws.Range("N32").Value = -15805.25
$ws.Range("H122").Value = 2136.4285
$ws.Range("I122").Value = 2086.8
$ws.Range("J122").Value = 2260.5
$ws.Range("K122").Value = 6260.400000000001
$ws.Range("L122").Value = 6781.5
$ws.Range("M122").Value = -3810.400000000001
$ws.Range("N122").Value = -11681.5
$ws.Range("H132").Value = 2296
$ws.Range("I132").Value = 1307.6
$ws.Range("J132").Value = 4767
$ws.Range("K132").Value = 3922.8
$ws.Range("L132").Value = 14301
$ws.Range("M132").Value = -1392.8
$ws.Range("N132").Value = -19361

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1759.1538
$ws.Range("I86").Value = 1641.3334
$ws.Range("J86").Value = 2024.25
$ws.Range("K86").Value = 1641.3334
$ws.Range("L86").Value = 2024.25
$ws.Range("M86").Value = -518.3334
$ws.Range("N86").Value = -4270.25
$ws.Range("H89").Value = 1759.1538
$ws.Range("I89").Value = 1641.3334
$ws.Range("J89").Value = 2024.25
$ws.Range("K89").Value = 8206.666999999999
$ws.Range("L89").Value = 10121.25
$ws.Range("M89").Value = -2590.666999999999
$ws.Range("N89").Value = -21353.25
$ws.Range("H94").Value = 1463.3334
$ws.Range("I94").Value = 851.4286
$ws.Range("J94").Value = 2320
$ws.Range("K94").Value = 851.4286
$ws.Range("L94").Value = 2320
$ws.Range("M94").Value = -400.4286
$ws.Range("N94").Value = -3222
$ws.Range("H105").Value = 1139508.1
$ws.Range("I105").Value = 1895788.5
$ws.Range("J105").Value = 5087.5
$ws.Range("K105").Value = 1895788.5
$ws.Range("L105").Value = 5087.5
$ws.Range("M105").Value = -1894041.5
$ws.Range("N105").Value = -8581.5
$ws.Range("H107").Value = 871.9524
$ws.Range("I107").Value = 854.7646999999999
$ws.Range("J107").Value = 945
$ws.Range("K107").Value = 854.7646999999999
$ws.Range("L107").Value = 945
$ws.Range("M107").Value = 1065.2353
$ws.Range("N107").Value = -4785
$ws.Range("H134").Value = 6251873
$ws.Range("I134").Value = 10871098
$ws.Range("J134").Value = 2333.8823
$ws.Range("K134").Value = 32613294
$ws.Range("L134").Value = 7001.646900000001
$ws.Range("M134").Value = -32610759
$ws.Range("N134").Value = -12071.6469

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9807162
$ws.Range("I31").Value = 21740300
$ws.Range("J31").Value = 4940.4287
$ws.Range("K31").Value = 21740300
$ws.Range("L31").Value = 4940.4287
$ws.Range("M31").Value = -21740005
$ws.Range("N31").Value = -5530.4287
$ws.Range("H34").Value = 9807162
$ws.Range("I34").Value = 21740300
$ws.Range("J34").Value = 4940.4287
$ws.Range("K34").Value = 21740300
$ws.Range("L34").Value = 4940.4287
$ws.Range("M34").Value = -21740098
$ws.Range("N34").Value = -5344.4287
$ws.Range("H107").Value = 1161.9259
$ws.Range("I107").Value = 650.8570999999999
$ws.Range("K107").Value = 650.8570999999999
$ws.Range("M107").Value = 1269.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 444.375
$ws.Range("I5").Value = 303.9091
$ws.Range("J5").Value = 753.4
$ws.Range("K5").Value = 911.7273
$ws.Range("L5").Value = 2260.2
$ws.Range("M5").Value = -799.7273
$ws.Range("N5").Value = -2484.2
$ws.Range("H115").Value = 1584.1538
$ws.Range("I115").Value = 1360.4849
$ws.Range("J115").Value = 2814.3333
$ws.Range("K115").Value = 4081.4547
$ws.Range("L115").Value = 8442.999899999999
$ws.Range("M115").Value = -2906.4547
$ws.Range("N115").Value = -10792.9999
$ws.Range("H120").Value = 14664.333
$ws.Range("I120").Value = 8250
$ws.Range("J120").Value = 19795.8
$ws.Range("K120").Value = 24750
$ws.Range("L120").Value = 59387.39999999999
$ws.Range("M120").Value = -19912
$ws.Range("N120").Value = -69063.39999999999
$ws.Range("H121").Value = 62507508
$ws.Range("I121").Value = 900
$ws.Range("J121").Value = 71437020
$ws.Range("K121").Value = 2700
$ws.Range("L121").Value = 214311060
$ws.Range("M121").Value = -1390
$ws.Range("N121").Value = -214313680
$ws.Range("H122").Value = 883.5217
$ws.Range("I122").Value = 730.6429000000001
$ws.Range("J122").Value = 1121.3334
$ws.Range("K122").Value = 6575.7861
$ws.Range("L122").Value = 10092.0006
$ws.Range("M122").Value = -4125.7861
$ws.Range("N122").Value = -14992.0006
$ws.Range("H135").Value = 444.375
$ws.Range("I135").Value = 303.9091
$ws.Range("J135").Value = 753.4
$ws.Range("K135").Value = 2735.1819
$ws.Range("L135").Value = 6780.599999999999
$ws.Range("M135").Value = -200.1819
$ws.Range("N135").Value = -11850.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 76437910
$ws.Range("J11").Value = 110006870
$ws.Range("L11").Value = 110006870
$ws.Range("N11").Value = -110007148
$ws.Range("H24").Value = 235670.94
$ws.Range("I24").Value = 4000006
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 4000006
$ws.Range("L24").Value = 400
$ws.Range("M24").Value = -3999833
$ws.Range("N24").Value = -746
$ws.Range("H80").Value = 2758.75
$ws.Range("I80").Value = 2630.5
$ws.Range("J80").Value = 3400
$ws.Range("K80").Value = 2630.5
$ws.Range("L80").Value = 3400
$ws.Range("M80").Value = -1632.5
$ws.Range("N80").Value = -5396
$ws.Range("H83").Value = 2758.75
$ws.Range("I83").Value = 2630.5
$ws.Range("J83").Value = 3400
$ws.Range("K83").Value = 13152.5
$ws.Range("L83").Value = 17000
$ws.Range("M83").Value = -8160.5
$ws.Range("N83").Value = -26984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 5933.3335
$ws.Range("J24").Value = 5933.3335
$ws.Range("L24").Value = 5933.3335
$ws.Range("N24").Value = -6619.3335
$ws.Range("H93").Value = 1834.7742
$ws.Range("I93").Value = 1663.0476
$ws.Range("J93").Value = 2195.4
$ws.Range("K93").Value = 1663.0476
$ws.Range("L93").Value = 2195.4
$ws.Range("M93").Value = -415.0476000000001
$ws.Range("N93").Value = -4691.4
$ws.Range("H136").Value = 2450.5
$ws.Range("I136").Value = 1213.5217
$ws.Range("J136").Value = 5611.6665
$ws.Range("K136").Value = 3640.5651
$ws.Range("L136").Value = 16834.9995
$ws.Range("M136").Value = -1090.5651
$ws.Range("N136").Value = -21934.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 588.3077
$ws.Range("I107").Value = 448.85715
$ws.Range("J107").Value = 751
$ws.Range("K107").Value = 1346.57145
$ws.Range("L107").Value = 2253
$ws.Range("M107").Value = 573.4285500000001
$ws.Range("N107").Value = -6093
